# Fix typo in presentation:
#   "[Xintao et al, 2020}: " -> "[Xiantao et al, 2020}: "
# on the "Game plan" slide (slide 13), in the bullet about further
# evaluation metrics (Signal-to-noise ratio / peak SNR).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$para = $tr.Paragraphs(11)
$run = $para.Runs(1)
$run.Text = "Add further evaluation metrics following [Xiantao et al, 2020}: "
